$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.287.53'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.411.28'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '560.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.28%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.55%  '
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("E10").Value = '  -0.89%  '
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.838.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.196.69'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.411.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '324.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.78%  '
$ws.Range("E21").Value = '  +1.82%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.76%  '
$ws.Range("E24").Value = '  +1.29%  '
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.14%  '
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.59'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.402'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.38'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.64%  '
$ws.Range("E35").Value = '  +3.52%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '323.70'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.07%  '
$ws.Range("E40").Value = '  -0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.70'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '147.80'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.07%  '
$ws.Range("E43").Value = '  -3.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0971'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("E45").Value = '  +1.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0517'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.574'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("E48").Value = '  -1.38%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("E51").Value = '  -0.67%  '
